## Sample Project / Main.xlsx - "Rules" sheet
## Change cell B11 from "R40" to "1".
##
## B11 keeps its existing (text) cell type and existing cell style (s="23")
## -- only the stored string changes, from the shared string "R40" to a
## brand-new shared string "1". Because "1" looks like a number, a plain
## `Range.Value = "1"` assignment would make Excel re-interpret it as a
## numeric literal (and typically also reshuffle the cell's number format /
## style). To avoid that, we enter it as a formula that evaluates to the
## text "1", then convert the cell in place to its evaluated value via
## Copy + PasteSpecial(xlPasteValues) -- this is the standard "convert
## formula to literal value without touching formatting" move, so the
## cell's style/format stay exactly as they were.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")

# Type a formula whose result is the text string "1" (not the number 1).
$cell.Formula = '="1"'

# Freeze the formula result into a literal value in place, leaving the
# cell's existing formatting/style untouched.
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues

# Clean up the copy "marching ants" mode.
$excel.CutCopyMode = 0
